# Re-create the auto-generated "ListLabel" character styles (647-777) that
# LibreOffice/Word mint whenever new numbered/bulleted list paragraphs are
# typed into the document (HTML memorisation notes, JS functional programming,
# Redux notes). Mirrors styles.xml additions from the commit.
$d = $word.ActiveDocument

$s647 = $d.Styles.Add("ListLabel647", 2)
$s647.NameLocal = "ListLabel 647"
$s647.Font.NameBi = "OpenSymbol"
$s647.Font.Bold = $true
$s647.QuickStyle = $true

$s648 = $d.Styles.Add("ListLabel648", 2)
$s648.NameLocal = "ListLabel 648"
$s648.Font.NameBi = "OpenSymbol"
$s648.Font.Bold = $true
$s648.QuickStyle = $true

$s649 = $d.Styles.Add("ListLabel649", 2)
$s649.NameLocal = "ListLabel 649"
$s649.Font.NameBi = "OpenSymbol"
$s649.QuickStyle = $true

$s650 = $d.Styles.Add("ListLabel650", 2)
$s650.NameLocal = "ListLabel 650"
$s650.Font.NameBi = "OpenSymbol"
$s650.QuickStyle = $true

$s651 = $d.Styles.Add("ListLabel651", 2)
$s651.NameLocal = "ListLabel 651"
$s651.Font.NameBi = "OpenSymbol"
$s651.QuickStyle = $true

$s652 = $d.Styles.Add("ListLabel652", 2)
$s652.NameLocal = "ListLabel 652"
$s652.Font.NameBi = "OpenSymbol"
$s652.QuickStyle = $true

$s653 = $d.Styles.Add("ListLabel653", 2)
$s653.NameLocal = "ListLabel 653"
$s653.Font.NameBi = "OpenSymbol"
$s653.QuickStyle = $true

$s654 = $d.Styles.Add("ListLabel654", 2)
$s654.NameLocal = "ListLabel 654"
$s654.Font.NameBi = "OpenSymbol"
$s654.QuickStyle = $true

$s655 = $d.Styles.Add("ListLabel655", 2)
$s655.NameLocal = "ListLabel 655"
$s655.Font.NameBi = "OpenSymbol"
$s655.QuickStyle = $true

$s656 = $d.Styles.Add("ListLabel656", 2)
$s656.NameLocal = "ListLabel 656"
$s656.Font.NameBi = "OpenSymbol"
$s656.Font.Bold = $false
$s656.QuickStyle = $true

$s657 = $d.Styles.Add("ListLabel657", 2)
$s657.NameLocal = "ListLabel 657"
$s657.Font.NameBi = "OpenSymbol"
$s657.QuickStyle = $true

$s658 = $d.Styles.Add("ListLabel658", 2)
$s658.NameLocal = "ListLabel 658"
$s658.Font.NameBi = "OpenSymbol"
$s658.QuickStyle = $true

$s659 = $d.Styles.Add("ListLabel659", 2)
$s659.NameLocal = "ListLabel 659"
$s659.Font.NameBi = "OpenSymbol"
$s659.QuickStyle = $true

$s660 = $d.Styles.Add("ListLabel660", 2)
$s660.NameLocal = "ListLabel 660"
$s660.Font.NameBi = "OpenSymbol"
$s660.QuickStyle = $true

$s661 = $d.Styles.Add("ListLabel661", 2)
$s661.NameLocal = "ListLabel 661"
$s661.Font.NameBi = "OpenSymbol"
$s661.QuickStyle = $true

$s662 = $d.Styles.Add("ListLabel662", 2)
$s662.NameLocal = "ListLabel 662"
$s662.Font.NameBi = "OpenSymbol"
$s662.QuickStyle = $true

$s663 = $d.Styles.Add("ListLabel663", 2)
$s663.NameLocal = "ListLabel 663"
$s663.Font.NameBi = "OpenSymbol"
$s663.QuickStyle = $true

$s664 = $d.Styles.Add("ListLabel664", 2)
$s664.NameLocal = "ListLabel 664"
$s664.Font.NameBi = "OpenSymbol"
$s664.QuickStyle = $true

$s665 = $d.Styles.Add("ListLabel665", 2)
$s665.NameLocal = "ListLabel 665"
$s665.Font.NameBi = "OpenSymbol"
$s665.Font.Bold = $false
$s665.QuickStyle = $true

$s666 = $d.Styles.Add("ListLabel666", 2)
$s666.NameLocal = "ListLabel 666"
$s666.Font.NameBi = "OpenSymbol"
$s666.QuickStyle = $true

$s667 = $d.Styles.Add("ListLabel667", 2)
$s667.NameLocal = "ListLabel 667"
$s667.Font.NameBi = "OpenSymbol"
$s667.QuickStyle = $true

$s668 = $d.Styles.Add("ListLabel668", 2)
$s668.NameLocal = "ListLabel 668"
$s668.Font.NameBi = "OpenSymbol"
$s668.QuickStyle = $true

$s669 = $d.Styles.Add("ListLabel669", 2)
$s669.NameLocal = "ListLabel 669"
$s669.Font.NameBi = "OpenSymbol"
$s669.QuickStyle = $true

$s670 = $d.Styles.Add("ListLabel670", 2)
$s670.NameLocal = "ListLabel 670"
$s670.Font.NameBi = "OpenSymbol"
$s670.QuickStyle = $true

$s671 = $d.Styles.Add("ListLabel671", 2)
$s671.NameLocal = "ListLabel 671"
$s671.Font.NameBi = "OpenSymbol"
$s671.QuickStyle = $true

$s672 = $d.Styles.Add("ListLabel672", 2)
$s672.NameLocal = "ListLabel 672"
$s672.Font.NameBi = "OpenSymbol"
$s672.QuickStyle = $true

$s673 = $d.Styles.Add("ListLabel673", 2)
$s673.NameLocal = "ListLabel 673"
$s673.Font.NameBi = "OpenSymbol"
$s673.QuickStyle = $true

$s674 = $d.Styles.Add("ListLabel674", 2)
$s674.NameLocal = "ListLabel 674"
$s674.Font.NameBi = "OpenSymbol"
$s674.Font.Bold = $false
$s674.QuickStyle = $true

$s675 = $d.Styles.Add("ListLabel675", 2)
$s675.NameLocal = "ListLabel 675"
$s675.Font.NameBi = "OpenSymbol"
$s675.QuickStyle = $true

$s676 = $d.Styles.Add("ListLabel676", 2)
$s676.NameLocal = "ListLabel 676"
$s676.Font.NameBi = "OpenSymbol"
$s676.QuickStyle = $true

$s677 = $d.Styles.Add("ListLabel677", 2)
$s677.NameLocal = "ListLabel 677"
$s677.Font.NameBi = "OpenSymbol"
$s677.QuickStyle = $true

$s678 = $d.Styles.Add("ListLabel678", 2)
$s678.NameLocal = "ListLabel 678"
$s678.Font.NameBi = "OpenSymbol"
$s678.QuickStyle = $true

$s679 = $d.Styles.Add("ListLabel679", 2)
$s679.NameLocal = "ListLabel 679"
$s679.Font.NameBi = "OpenSymbol"
$s679.QuickStyle = $true

$s680 = $d.Styles.Add("ListLabel680", 2)
$s680.NameLocal = "ListLabel 680"
$s680.Font.NameBi = "OpenSymbol"
$s680.QuickStyle = $true

$s681 = $d.Styles.Add("ListLabel681", 2)
$s681.NameLocal = "ListLabel 681"
$s681.Font.NameBi = "OpenSymbol"
$s681.QuickStyle = $true

$s682 = $d.Styles.Add("ListLabel682", 2)
$s682.NameLocal = "ListLabel 682"
$s682.Font.NameBi = "OpenSymbol"
$s682.QuickStyle = $true

$s683 = $d.Styles.Add("ListLabel683", 2)
$s683.NameLocal = "ListLabel 683"
$s683.Font.NameBi = "OpenSymbol"
$s683.Font.Bold = $false
$s683.QuickStyle = $true

$s684 = $d.Styles.Add("ListLabel684", 2)
$s684.NameLocal = "ListLabel 684"
$s684.Font.NameBi = "OpenSymbol"
$s684.QuickStyle = $true

$s685 = $d.Styles.Add("ListLabel685", 2)
$s685.NameLocal = "ListLabel 685"
$s685.Font.NameBi = "OpenSymbol"
$s685.QuickStyle = $true

$s686 = $d.Styles.Add("ListLabel686", 2)
$s686.NameLocal = "ListLabel 686"
$s686.Font.NameBi = "OpenSymbol"
$s686.QuickStyle = $true

$s687 = $d.Styles.Add("ListLabel687", 2)
$s687.NameLocal = "ListLabel 687"
$s687.Font.NameBi = "OpenSymbol"
$s687.QuickStyle = $true

$s688 = $d.Styles.Add("ListLabel688", 2)
$s688.NameLocal = "ListLabel 688"
$s688.Font.NameBi = "OpenSymbol"
$s688.QuickStyle = $true

$s689 = $d.Styles.Add("ListLabel689", 2)
$s689.NameLocal = "ListLabel 689"
$s689.Font.NameBi = "OpenSymbol"
$s689.QuickStyle = $true

$s690 = $d.Styles.Add("ListLabel690", 2)
$s690.NameLocal = "ListLabel 690"
$s690.Font.NameBi = "OpenSymbol"
$s690.QuickStyle = $true

$s691 = $d.Styles.Add("ListLabel691", 2)
$s691.NameLocal = "ListLabel 691"
$s691.Font.NameBi = "OpenSymbol"
$s691.QuickStyle = $true

$s692 = $d.Styles.Add("ListLabel692", 2)
$s692.NameLocal = "ListLabel 692"
$s692.Font.NameBi = "OpenSymbol"
$s692.Font.Bold = $false
$s692.QuickStyle = $true

$s693 = $d.Styles.Add("ListLabel693", 2)
$s693.NameLocal = "ListLabel 693"
$s693.Font.NameBi = "OpenSymbol"
$s693.QuickStyle = $true

$s694 = $d.Styles.Add("ListLabel694", 2)
$s694.NameLocal = "ListLabel 694"
$s694.Font.NameBi = "OpenSymbol"
$s694.QuickStyle = $true

$s695 = $d.Styles.Add("ListLabel695", 2)
$s695.NameLocal = "ListLabel 695"
$s695.Font.NameBi = "OpenSymbol"
$s695.QuickStyle = $true

$s696 = $d.Styles.Add("ListLabel696", 2)
$s696.NameLocal = "ListLabel 696"
$s696.Font.NameBi = "OpenSymbol"
$s696.QuickStyle = $true

$s697 = $d.Styles.Add("ListLabel697", 2)
$s697.NameLocal = "ListLabel 697"
$s697.Font.NameBi = "OpenSymbol"
$s697.QuickStyle = $true

$s698 = $d.Styles.Add("ListLabel698", 2)
$s698.NameLocal = "ListLabel 698"
$s698.Font.NameBi = "OpenSymbol"
$s698.QuickStyle = $true

$s699 = $d.Styles.Add("ListLabel699", 2)
$s699.NameLocal = "ListLabel 699"
$s699.Font.NameBi = "OpenSymbol"
$s699.QuickStyle = $true

$s700 = $d.Styles.Add("ListLabel700", 2)
$s700.NameLocal = "ListLabel 700"
$s700.Font.NameBi = "OpenSymbol"
$s700.QuickStyle = $true

$s701 = $d.Styles.Add("ListLabel701", 2)
$s701.NameLocal = "ListLabel 701"
$s701.Font.NameBi = "OpenSymbol"
$s701.Font.Bold = $false
$s701.QuickStyle = $true

$s702 = $d.Styles.Add("ListLabel702", 2)
$s702.NameLocal = "ListLabel 702"
$s702.Font.NameBi = "OpenSymbol"
$s702.QuickStyle = $true

$s703 = $d.Styles.Add("ListLabel703", 2)
$s703.NameLocal = "ListLabel 703"
$s703.Font.NameBi = "OpenSymbol"
$s703.QuickStyle = $true

$s704 = $d.Styles.Add("ListLabel704", 2)
$s704.NameLocal = "ListLabel 704"
$s704.Font.NameBi = "OpenSymbol"
$s704.QuickStyle = $true

$s705 = $d.Styles.Add("ListLabel705", 2)
$s705.NameLocal = "ListLabel 705"
$s705.Font.NameBi = "OpenSymbol"
$s705.QuickStyle = $true

$s706 = $d.Styles.Add("ListLabel706", 2)
$s706.NameLocal = "ListLabel 706"
$s706.Font.NameBi = "OpenSymbol"
$s706.QuickStyle = $true

$s707 = $d.Styles.Add("ListLabel707", 2)
$s707.NameLocal = "ListLabel 707"
$s707.Font.NameBi = "OpenSymbol"
$s707.QuickStyle = $true

$s708 = $d.Styles.Add("ListLabel708", 2)
$s708.NameLocal = "ListLabel 708"
$s708.Font.NameBi = "OpenSymbol"
$s708.QuickStyle = $true

$s709 = $d.Styles.Add("ListLabel709", 2)
$s709.NameLocal = "ListLabel 709"
$s709.Font.NameBi = "OpenSymbol"
$s709.QuickStyle = $true

$s710 = $d.Styles.Add("ListLabel710", 2)
$s710.NameLocal = "ListLabel 710"
$s710.Font.NameBi = "OpenSymbol"
$s710.Font.Bold = $false
$s710.QuickStyle = $true

$s711 = $d.Styles.Add("ListLabel711", 2)
$s711.NameLocal = "ListLabel 711"
$s711.Font.NameBi = "OpenSymbol"
$s711.QuickStyle = $true

$s712 = $d.Styles.Add("ListLabel712", 2)
$s712.NameLocal = "ListLabel 712"
$s712.Font.NameBi = "OpenSymbol"
$s712.QuickStyle = $true

$s713 = $d.Styles.Add("ListLabel713", 2)
$s713.NameLocal = "ListLabel 713"
$s713.Font.NameBi = "OpenSymbol"
$s713.QuickStyle = $true

$s714 = $d.Styles.Add("ListLabel714", 2)
$s714.NameLocal = "ListLabel 714"
$s714.Font.NameBi = "OpenSymbol"
$s714.QuickStyle = $true

$s715 = $d.Styles.Add("ListLabel715", 2)
$s715.NameLocal = "ListLabel 715"
$s715.Font.NameBi = "OpenSymbol"
$s715.QuickStyle = $true

$s716 = $d.Styles.Add("ListLabel716", 2)
$s716.NameLocal = "ListLabel 716"
$s716.Font.NameBi = "OpenSymbol"
$s716.QuickStyle = $true

$s717 = $d.Styles.Add("ListLabel717", 2)
$s717.NameLocal = "ListLabel 717"
$s717.Font.NameBi = "OpenSymbol"
$s717.QuickStyle = $true

$s718 = $d.Styles.Add("ListLabel718", 2)
$s718.NameLocal = "ListLabel 718"
$s718.Font.NameBi = "OpenSymbol"
$s718.QuickStyle = $true

$s719 = $d.Styles.Add("ListLabel719", 2)
$s719.NameLocal = "ListLabel 719"
$s719.Font.NameBi = "OpenSymbol"
$s719.Font.Bold = $false
$s719.QuickStyle = $true

$s720 = $d.Styles.Add("ListLabel720", 2)
$s720.NameLocal = "ListLabel 720"
$s720.Font.NameBi = "OpenSymbol"
$s720.QuickStyle = $true

$s721 = $d.Styles.Add("ListLabel721", 2)
$s721.NameLocal = "ListLabel 721"
$s721.Font.NameBi = "OpenSymbol"
$s721.QuickStyle = $true

$s722 = $d.Styles.Add("ListLabel722", 2)
$s722.NameLocal = "ListLabel 722"
$s722.Font.NameBi = "OpenSymbol"
$s722.QuickStyle = $true

$s723 = $d.Styles.Add("ListLabel723", 2)
$s723.NameLocal = "ListLabel 723"
$s723.Font.NameBi = "OpenSymbol"
$s723.QuickStyle = $true

$s724 = $d.Styles.Add("ListLabel724", 2)
$s724.NameLocal = "ListLabel 724"
$s724.Font.NameBi = "OpenSymbol"
$s724.QuickStyle = $true

$s725 = $d.Styles.Add("ListLabel725", 2)
$s725.NameLocal = "ListLabel 725"
$s725.Font.NameBi = "OpenSymbol"
$s725.QuickStyle = $true

$s726 = $d.Styles.Add("ListLabel726", 2)
$s726.NameLocal = "ListLabel 726"
$s726.Font.NameBi = "OpenSymbol"
$s726.QuickStyle = $true

$s727 = $d.Styles.Add("ListLabel727", 2)
$s727.NameLocal = "ListLabel 727"
$s727.Font.NameBi = "OpenSymbol"
$s727.QuickStyle = $true

$s728 = $d.Styles.Add("ListLabel728", 2)
$s728.NameLocal = "ListLabel 728"
$s728.Font.NameBi = "OpenSymbol"
$s728.Font.Bold = $false
$s728.QuickStyle = $true

$s729 = $d.Styles.Add("ListLabel729", 2)
$s729.NameLocal = "ListLabel 729"
$s729.Font.NameBi = "OpenSymbol"
$s729.QuickStyle = $true

$s730 = $d.Styles.Add("ListLabel730", 2)
$s730.NameLocal = "ListLabel 730"
$s730.Font.NameBi = "OpenSymbol"
$s730.QuickStyle = $true

$s731 = $d.Styles.Add("ListLabel731", 2)
$s731.NameLocal = "ListLabel 731"
$s731.Font.NameBi = "OpenSymbol"
$s731.QuickStyle = $true

$s732 = $d.Styles.Add("ListLabel732", 2)
$s732.NameLocal = "ListLabel 732"
$s732.Font.NameBi = "OpenSymbol"
$s732.QuickStyle = $true

$s733 = $d.Styles.Add("ListLabel733", 2)
$s733.NameLocal = "ListLabel 733"
$s733.Font.NameBi = "OpenSymbol"
$s733.QuickStyle = $true

$s734 = $d.Styles.Add("ListLabel734", 2)
$s734.NameLocal = "ListLabel 734"
$s734.Font.NameBi = "OpenSymbol"
$s734.QuickStyle = $true

$s735 = $d.Styles.Add("ListLabel735", 2)
$s735.NameLocal = "ListLabel 735"
$s735.Font.NameBi = "OpenSymbol"
$s735.QuickStyle = $true

$s736 = $d.Styles.Add("ListLabel736", 2)
$s736.NameLocal = "ListLabel 736"
$s736.Font.NameBi = "OpenSymbol"
$s736.QuickStyle = $true

$s737 = $d.Styles.Add("ListLabel737", 2)
$s737.NameLocal = "ListLabel 737"
$s737.Font.NameBi = "OpenSymbol"
$s737.Font.Bold = $false
$s737.QuickStyle = $true

$s738 = $d.Styles.Add("ListLabel738", 2)
$s738.NameLocal = "ListLabel 738"
$s738.Font.NameBi = "OpenSymbol"
$s738.QuickStyle = $true

$s739 = $d.Styles.Add("ListLabel739", 2)
$s739.NameLocal = "ListLabel 739"
$s739.Font.NameBi = "OpenSymbol"
$s739.QuickStyle = $true

$s740 = $d.Styles.Add("ListLabel740", 2)
$s740.NameLocal = "ListLabel 740"
$s740.Font.NameBi = "OpenSymbol"
$s740.QuickStyle = $true

$s741 = $d.Styles.Add("ListLabel741", 2)
$s741.NameLocal = "ListLabel 741"
$s741.Font.NameBi = "OpenSymbol"
$s741.QuickStyle = $true

$s742 = $d.Styles.Add("ListLabel742", 2)
$s742.NameLocal = "ListLabel 742"
$s742.Font.NameBi = "OpenSymbol"
$s742.QuickStyle = $true

$s743 = $d.Styles.Add("ListLabel743", 2)
$s743.NameLocal = "ListLabel 743"
$s743.Font.NameBi = "OpenSymbol"
$s743.QuickStyle = $true

$s744 = $d.Styles.Add("ListLabel744", 2)
$s744.NameLocal = "ListLabel 744"
$s744.Font.NameBi = "OpenSymbol"
$s744.QuickStyle = $true

$s745 = $d.Styles.Add("ListLabel745", 2)
$s745.NameLocal = "ListLabel 745"
$s745.Font.NameBi = "OpenSymbol"
$s745.QuickStyle = $true

$s746 = $d.Styles.Add("ListLabel746", 2)
$s746.NameLocal = "ListLabel 746"
$s746.Font.NameBi = "OpenSymbol"
$s746.QuickStyle = $true

$s747 = $d.Styles.Add("ListLabel747", 2)
$s747.NameLocal = "ListLabel 747"
$s747.Font.NameBi = "OpenSymbol"
$s747.Font.Bold = $false
$s747.QuickStyle = $true

$s748 = $d.Styles.Add("ListLabel748", 2)
$s748.NameLocal = "ListLabel 748"
$s748.Font.NameBi = "OpenSymbol"
$s748.QuickStyle = $true

$s749 = $d.Styles.Add("ListLabel749", 2)
$s749.NameLocal = "ListLabel 749"
$s749.Font.NameBi = "OpenSymbol"
$s749.QuickStyle = $true

$s750 = $d.Styles.Add("ListLabel750", 2)
$s750.NameLocal = "ListLabel 750"
$s750.Font.NameBi = "OpenSymbol"
$s750.QuickStyle = $true

$s751 = $d.Styles.Add("ListLabel751", 2)
$s751.NameLocal = "ListLabel 751"
$s751.Font.NameBi = "OpenSymbol"
$s751.QuickStyle = $true

$s752 = $d.Styles.Add("ListLabel752", 2)
$s752.NameLocal = "ListLabel 752"
$s752.Font.NameBi = "OpenSymbol"
$s752.QuickStyle = $true

$s753 = $d.Styles.Add("ListLabel753", 2)
$s753.NameLocal = "ListLabel 753"
$s753.Font.NameBi = "OpenSymbol"
$s753.QuickStyle = $true

$s754 = $d.Styles.Add("ListLabel754", 2)
$s754.NameLocal = "ListLabel 754"
$s754.Font.NameBi = "OpenSymbol"
$s754.QuickStyle = $true

$s755 = $d.Styles.Add("ListLabel755", 2)
$s755.NameLocal = "ListLabel 755"
$s755.Font.NameBi = "OpenSymbol"
$s755.Font.Bold = $false
$s755.QuickStyle = $true

$s756 = $d.Styles.Add("ListLabel756", 2)
$s756.NameLocal = "ListLabel 756"
$s756.Font.NameBi = "OpenSymbol"
$s756.Font.Bold = $false
$s756.QuickStyle = $true

$s757 = $d.Styles.Add("ListLabel757", 2)
$s757.NameLocal = "ListLabel 757"
$s757.Font.NameBi = "OpenSymbol"
$s757.QuickStyle = $true

$s758 = $d.Styles.Add("ListLabel758", 2)
$s758.NameLocal = "ListLabel 758"
$s758.Font.NameBi = "OpenSymbol"
$s758.QuickStyle = $true

$s759 = $d.Styles.Add("ListLabel759", 2)
$s759.NameLocal = "ListLabel 759"
$s759.Font.NameBi = "OpenSymbol"
$s759.QuickStyle = $true

$s760 = $d.Styles.Add("ListLabel760", 2)
$s760.NameLocal = "ListLabel 760"
$s760.Font.NameBi = "OpenSymbol"
$s760.QuickStyle = $true

$s761 = $d.Styles.Add("ListLabel761", 2)
$s761.NameLocal = "ListLabel 761"
$s761.Font.NameBi = "OpenSymbol"
$s761.QuickStyle = $true

$s762 = $d.Styles.Add("ListLabel762", 2)
$s762.NameLocal = "ListLabel 762"
$s762.Font.NameBi = "OpenSymbol"
$s762.QuickStyle = $true

$s763 = $d.Styles.Add("ListLabel763", 2)
$s763.NameLocal = "ListLabel 763"
$s763.Font.NameBi = "OpenSymbol"
$s763.QuickStyle = $true

$s764 = $d.Styles.Add("ListLabel764", 2)
$s764.NameLocal = "ListLabel 764"
$s764.Font.NameBi = "OpenSymbol"
$s764.QuickStyle = $true

$s765 = $d.Styles.Add("ListLabel765", 2)
$s765.NameLocal = "ListLabel 765"
$s765.Font.NameBi = "OpenSymbol"
$s765.QuickStyle = $true

$s766 = $d.Styles.Add("ListLabel766", 2)
$s766.NameLocal = "ListLabel 766"
$s766.Font.NameBi = "OpenSymbol"
$s766.QuickStyle = $true

$s767 = $d.Styles.Add("ListLabel767", 2)
$s767.NameLocal = "ListLabel 767"
$s767.Font.NameBi = "OpenSymbol"
$s767.QuickStyle = $true

$s768 = $d.Styles.Add("ListLabel768", 2)
$s768.NameLocal = "ListLabel 768"
$s768.Font.NameBi = "OpenSymbol"
$s768.QuickStyle = $true

$s769 = $d.Styles.Add("ListLabel769", 2)
$s769.NameLocal = "ListLabel 769"
$s769.Font.NameBi = "OpenSymbol"
$s769.QuickStyle = $true

$s770 = $d.Styles.Add("ListLabel770", 2)
$s770.NameLocal = "ListLabel 770"
$s770.Font.NameBi = "OpenSymbol"
$s770.QuickStyle = $true

$s771 = $d.Styles.Add("ListLabel771", 2)
$s771.NameLocal = "ListLabel 771"
$s771.Font.NameBi = "OpenSymbol"
$s771.QuickStyle = $true

$s772 = $d.Styles.Add("ListLabel772", 2)
$s772.NameLocal = "ListLabel 772"
$s772.Font.NameBi = "OpenSymbol"
$s772.QuickStyle = $true

$s773 = $d.Styles.Add("ListLabel773", 2)
$s773.NameLocal = "ListLabel 773"
$s773.Font.Bold = $false
$s773.Font.BoldBi = $false
$s773.Font.Color = 655360
$s773.QuickStyle = $true

$s774 = $d.Styles.Add("ListLabel774", 2)
$s774.NameLocal = "ListLabel 774"
$s774.QuickStyle = $true

$s775 = $d.Styles.Add("ListLabel775", 2)
$s775.NameLocal = "ListLabel 775"
$s775.Font.Bold = $false
$s775.Font.BoldBi = $false
$s775.Font.Color = 655360
$s775.QuickStyle = $true

$s776 = $d.Styles.Add("ListLabel776", 2)
$s776.NameLocal = "ListLabel 776"
$s776.QuickStyle = $true

$s777 = $d.Styles.Add("ListLabel777", 2)
$s777.NameLocal = "ListLabel 777"
$s777.Font.Color = 655360
$s777.QuickStyle = $true

Write-Output "Added $($d.Styles.Count) styles total (ListLabel647-ListLabel777 inserted)."
